$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sat right after the
#    "... dang nhap" heading paragraph (it gets relocated later, see
#    step 3, to sit inside the "password" row's "Textfield" cell).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# Helper: split the "Textfied" cell text (table 2 = the dang-nhap
# description table) into "Textfie" + "l" + "d", mirroring a user
# placing the caret after "Textfie" and typing the missing "l".
# Returns the Range for the freshly inserted "l" run so a caller can
# (optionally) drop a bookmark right after it.
# ------------------------------------------------------------------
function Split-Textfied($cell) {
    $cellRange = $cell.Range
    $cellRange.Text = "Textfie"

    $afterFie = $cell.Range
    $insertPos = $d.Range($afterFie.End - 1, $afterFie.End - 1)
    $insertPos.InsertAfter("l")
    # Force the new "l" text into its own run (distinct <w:r>) instead
    # of being silently merged back into the "Textfie" run.
    $insertPos.Bold = 1
    $insertPos.Bold = 0

    $lRange = $insertPos

    $afterL = $cell.Range
    $insertPos2 = $d.Range($afterL.End - 1, $afterL.End - 1)
    $insertPos2.InsertAfter("d")
    $insertPos2.Bold = 1
    $insertPos2.Bold = 0

    return $lRange
}

# Table 2 (1-based) is the "dang nhap" field/description table.
$t = $d.Tables.Item(2)

# ------------------------------------------------------------------
# 2. Row 2 (username row) - cell 3 holds "Textfied" -> "Textfield".
# ------------------------------------------------------------------
$cellUser = $t.Rows.Item(2).Cells.Item(3)
Split-Textfied $cellUser | Out-Null

# ------------------------------------------------------------------
# 3. Row 3 (password row) - cell 3 holds "Textfied" -> "Textfield",
#    and this time the editor's last caret position (the "_GoBack"
#    bookmark) lands right after the inserted "l".
# ------------------------------------------------------------------
$cellPass = $t.Rows.Item(3).Cells.Item(3)
$lRange = Split-Textfied $cellPass

$d.Bookmarks.Add("_GoBack", $lRange) | Out-Null

# ------------------------------------------------------------------
# 4. Row 4 (checkbox row) tblPrEx gains a tblCellMar override
#    (matching the rest of the table's per-row overrides).
# ------------------------------------------------------------------
$row4 = $t.Rows.Item(4)
$row4.LeftPadding = 108
$row4.RightPadding = 108
$row4.TopPadding = 0
$row4.BottomPadding = 0
